$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "timestamp" column (Z) values to reflect the latest run.
$ws.Range("Z2:Z22").Value = "2025-11-02T02:04:37.338374"
$ws.Range("Z23:Z26").Value = "2025-11-02T02:04:37.345485"
$ws.Range("Z27:Z29").Value = "2025-11-02T02:04:37.346477"
